$d = $word.ActiveDocument

function Insert-ParagraphXml($paraIndex, $pPrXml, $bodyXml, $plainText) {
    $p = $d.Paragraphs($paraIndex)
    $startPos = $p.Range.Start
    $insertPoint = $d.Range($startPos, $startPos)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrXml + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml)
    $deleteStart = $startPos + $plainText.Length
    $pAfter = $d.Paragraphs($paraIndex)
    $deleteEnd = $pAfter.Range.End - 1
    if ($deleteEnd -gt $deleteStart) {
        $toDelete = $d.Range($deleteStart, $deleteEnd)
        $toDelete.Delete()
    }
}

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# --- 1. Remove the "_GoBack" bookmark currently after "User Stories" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Add a trailing space run after "I'm able to play the game" (story 1) ---
$rngFind = $d.Content
$rngFind.Find.Execute("I" + [char]0x2019 + "m able to play the game", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngFind.Collapse(0)
$rngFind.InsertAfter(" ")

# --- 3. Story 3 (grid re-sizing, with grammar-check markers) ---
$body3 = '<w:r><w:t xml:space="preserve">As a developer, I want to have a re-sizeable grid so that more players </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>are able to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> play the game</w:t></w:r>'
$text3 = "As a developer, I want to have a re-sizeable grid so that more players are able to play the game"
Insert-ParagraphXml 3 $listPPr $body3 $text3

# --- 4. Story 4 (UDP server) ---
$body4 = '<w:r><w:t xml:space="preserve">As a </w:t></w:r><w:r><w:t>developer,</w:t></w:r><w:r><w:t xml:space="preserve"> I want </w:t></w:r><w:r><w:t>to implement a UDP server so that there can be multiple players</w:t></w:r>'
$text4 = "As a developer, I want to implement a UDP server so that there can be multiple players"
Insert-ParagraphXml 4 $listPPr $body4 $text4

# --- 5. Story 5 (leader boards, with a re-inserted _GoBack bookmark mid-sentence) ---
$body5 = '<w:r><w:t>As a developer</w:t></w:r><w:r><w:t xml:space="preserve">, I want </w:t></w:r><w:r><w:t>to imp</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>lement leader boards so that payers can see how they rank</w:t></w:r>'
$text5 = "As a developer, I want to implement leader boards so that payers can see how they rank"
Insert-ParagraphXml 5 $listPPr $body5 $text5
